$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.408218
$ws.Range("H2").Value = 61.224654
$ws.Range("I2").Value = 0.1108535210972707
$ws.Range("J2").Value = 0.1108535210972707
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4702473333333333
$ws.Range("N2").Value = 1.410742
$ws.Range("O2").Value = 0.009034922268422819
$ws.Range("P2").Value = 0.009034922268422819
$ws.Range("Q2").Value = 9.596910092585333
$ws.Range("R2").Value = 86.37219083326799
$ws.Range("S2").Value = 0.00100155294629481
$ws.Range("T2").Value = 0.00100155294629481

$ws.Range("G3").Value = 20.408218
$ws.Range("H3").Value = 61.224654
$ws.Range("I3").Value = 0.1108535210972707
$ws.Range("J3").Value = 0.1108535210972707
$ws.Range("N3").Value = 0.9584440000000001
$ws.Range("O3").Value = 0.006138235792679485
$ws.Range("P3").Value = 0.006138235792679485
$ws.Range("Q3").Value = 6.520044697597334
$ws.Range("R3").Value = 58.68040227837601
$ws.Range("S3").Value = 0.0006804450509438174
$ws.Range("T3").Value = 0.0006804450509438174

$ws.Range("G4").Value = 20.408218
$ws.Range("H4").Value = 61.224654
$ws.Range("I4").Value = 0.1108535210972707
$ws.Range("J4").Value = 0.1108535210972707
$ws.Range("M4").Value = 1.047307
$ws.Range("N4").Value = 3.141921
$ws.Range("O4").Value = 0.02012204358311108
$ws.Range("P4").Value = 0.02012204358311108
$ws.Range("Q4").Value = 21.373669568926
$ws.Range("R4").Value = 192.363026120334
$ws.Range("S4").Value = 0.002230599382860605
$ws.Range("T4").Value = 0.002230599382860605

$ws.Range("G5").Value = 20.408218
$ws.Range("H5").Value = 61.224654
$ws.Range("I5").Value = 0.1108535210972707
$ws.Range("J5").Value = 0.1108535210972707
$ws.Range("M5").Value = 50.21070966666667
$ws.Range("N5").Value = 150.632129
$ws.Range("O5").Value = 0.9647047983557866
$ws.Range("P5").Value = 0.9647047983557866
$ws.Range("Q5").Value = 1024.711108812041
$ws.Range("R5").Value = 9222.399979308368
$ws.Range("S5").Value = 0.1069409237171715
$ws.Range("T5").Value = 0.1069409237171715

$ws.Range("I6").Value = 0.2566851044076959
$ws.Range("J6").Value = 0.256685104407696
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4702473333333333
$ws.Range("N6").Value = 1.410742
$ws.Range("O6").Value = 0.009034922268422819
$ws.Range("P6").Value = 0.009034922268422819
$ws.Range("Q6").Value = 22.22197224520266
$ws.Range("R6").Value = 199.997750206824
$ws.Range("S6").Value = 0.002319129965785528
$ws.Range("T6").Value = 0.002319129965785529

$ws.Range("I7").Value = 0.2566851044076959
$ws.Range("J7").Value = 0.256685104407696
$ws.Range("N7").Value = 0.9584440000000001
$ws.Range("O7").Value = 0.006138235792679485
$ws.Range("P7").Value = 0.006138235792679485
$ws.Range("Q7").Value = 15.09738560741867
$ws.Range("S7").Value = 0.00157559369532299
$ws.Range("T7").Value = 0.00157559369532299

$ws.Range("I8").Value = 0.2566851044076959
$ws.Range("J8").Value = 0.256685104407696
$ws.Range("M8").Value = 1.047307
$ws.Range("N8").Value = 3.141921
$ws.Range("O8").Value = 0.02012204358311108
$ws.Range("P8").Value = 0.02012204358311108
$ws.Range("Q8").Value = 49.49145999666799
$ws.Range("R8").Value = 445.4231399700119
$ws.Range("S8").Value = 0.005165028858027076
$ws.Range("T8").Value = 0.005165028858027078

$ws.Range("I9").Value = 0.2566851044076959
$ws.Range("J9").Value = 0.256685104407696
$ws.Range("M9").Value = 50.21070966666667
$ws.Range("N9").Value = 150.632129
$ws.Range("O9").Value = 0.9647047983557866
$ws.Range("P9").Value = 0.9647047983557866
$ws.Range("Q9").Value = 2372.753479994065
$ws.Range("R9").Value = 21354.78131994659
$ws.Range("S9").Value = 0.2476253518885603
$ws.Range("T9").Value = 0.2476253518885604

$ws.Range("G10").Value = 85.307233
$ws.Range("H10").Value = 255.921699
$ws.Range("I10").Value = 0.4633725077375833
$ws.Range("J10").Value = 0.4633725077375833
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4702473333333333
$ws.Range("N10").Value = 1.410742
$ws.Range("O10").Value = 0.009034922268422819
$ws.Range("P10").Value = 0.009034922268422819
$ws.Range("Q10").Value = 40.11549883229533
$ws.Range("R10").Value = 361.039489490658
$ws.Range("S10").Value = 0.004186534588733216
$ws.Range("T10").Value = 0.004186534588733217

$ws.Range("G11").Value = 85.307233
$ws.Range("H11").Value = 255.921699
$ws.Range("I11").Value = 0.4633725077375833
$ws.Range("J11").Value = 0.4633725077375833
$ws.Range("N11").Value = 0.9584440000000001
$ws.Range("O11").Value = 0.006138235792679485
$ws.Range("P11").Value = 0.006138235792679485
$ws.Range("Q11").Value = 27.25406854181733
$ws.Range("R11").Value = 245.286616876356
$ws.Range("S11").Value = 0.002844289712338485
$ws.Range("T11").Value = 0.002844289712338486

$ws.Range("G12").Value = 85.307233
$ws.Range("H12").Value = 255.921699
$ws.Range("I12").Value = 0.4633725077375833
$ws.Range("J12").Value = 0.4633725077375833
$ws.Range("M12").Value = 1.047307
$ws.Range("N12").Value = 3.141921
$ws.Range("O12").Value = 0.02012204358311108
$ws.Range("P12").Value = 0.02012204358311108
$ws.Range("Q12").Value = 89.34286227153099
$ws.Range("R12").Value = 804.085760443779
$ws.Range("S12").Value = 0.009324001795911127
$ws.Range("T12").Value = 0.009324001795911129

$ws.Range("G13").Value = 85.307233
$ws.Range("H13").Value = 255.921699
$ws.Range("I13").Value = 0.4633725077375833
$ws.Range("J13").Value = 0.4633725077375833
$ws.Range("M13").Value = 50.21070966666667
$ws.Range("N13").Value = 150.632129
$ws.Range("O13").Value = 0.9647047983557866
$ws.Range("P13").Value = 0.9647047983557866
$ws.Range("Q13").Value = 4283.336708629686
$ws.Range("R13").Value = 38550.03037766717
$ws.Range("S13").Value = 0.4470176816406004
$ws.Range("T13").Value = 0.4470176816406005

$ws.Range("G14").Value = 31.12938966666666
$ws.Range("H14").Value = 93.38816899999999
$ws.Range("I14").Value = 0.16908886675745
$ws.Range("J14").Value = 0.16908886675745
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4702473333333333
$ws.Range("N14").Value = 1.410742
$ws.Range("O14").Value = 0.009034922268422819
$ws.Range("P14").Value = 0.009034922268422819
$ws.Range("Q14").Value = 14.63851247904422
$ws.Range("R14").Value = 131.746612311398
$ws.Range("S14").Value = 0.001527704767609264
$ws.Range("T14").Value = 0.001527704767609264

$ws.Range("G15").Value = 31.12938966666666
$ws.Range("H15").Value = 93.38816899999999
$ws.Range("I15").Value = 0.16908886675745
$ws.Range("J15").Value = 0.16908886675745
$ws.Range("N15").Value = 0.9584440000000001
$ws.Range("O15").Value = 0.006138235792679485
$ws.Range("P15").Value = 0.006138235792679485
$ws.Range("Q15").Value = 9.945258916559554
$ws.Range("R15").Value = 89.507330249036
$ws.Range("S15").Value = 0.001037907334074192
$ws.Range("T15").Value = 0.001037907334074192

$ws.Range("G16").Value = 31.12938966666666
$ws.Range("H16").Value = 93.38816899999999
$ws.Range("I16").Value = 0.16908886675745
$ws.Range("J16").Value = 0.16908886675745
$ws.Range("M16").Value = 1.047307
$ws.Range("N16").Value = 3.141921
$ws.Range("O16").Value = 0.02012204358311108
$ws.Range("P16").Value = 0.02012204358311108
$ws.Range("Q16").Value = 32.60202770362766
$ws.Range("R16").Value = 293.4182493326489
$ws.Range("S16").Value = 0.003402413546312272
$ws.Range("T16").Value = 0.003402413546312272

$ws.Range("G17").Value = 31.12938966666666
$ws.Range("H17").Value = 93.38816899999999
$ws.Range("I17").Value = 0.16908886675745
$ws.Range("J17").Value = 0.16908886675745
$ws.Range("M17").Value = 50.21070966666667
$ws.Range("N17").Value = 150.632129
$ws.Range("O17").Value = 0.9647047983557866
$ws.Range("P17").Value = 0.9647047983557866
$ws.Range("Q17").Value = 1563.028746653534
$ws.Range("R17").Value = 14067.2587198818
$ws.Range("S17").Value = 0.1631208411094543
$ws.Range("T17").Value = 0.1631208411094543
